$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 50

$ws.Cells.Item($row, 1).Value = "Record"
$ws.Cells.Item($row, 2).Value = "Balanço Geral"
$ws.Cells.Item($row, 3).Value = "Infraestrutura"
$ws.Cells.Item($row, 4).Value = "2025-04-04T13:38"
$ws.Cells.Item($row, 5).Value = "Negativo"
$ws.Cells.Item($row, 6).Value = 'Tô na bronca. Terrenos abandonados e crateras causam preocupação no bairro Alphaville. ´Rua Durval de Souza perto da esquina com Castro Alves. Muito mato e cratera na rua, tendo os moradores que sinalizar com galhos. Exibido vídeo enviado por morador, reclamando das condições. Segundo ele, moradores pediram para cortar galhos de árvores e, até o momento, não foram retirar. População também não colabora. Em outro ponto do bairro, única secretaria que atendeu foi a Posturas, mas disse que nada pode fazer. Acidentes estão acontecendo no local. Morador também procurou a Águas do Paraíba. Equipe pediu nota à prefeitura, mas não teve resposta ainda. '
